# Natmi following Dr Hou advice
# Update the Sostdc1-Lrp6 NATMI ligand-receptor pair table:
# the number of ligand-/receptor-expressing cells changed from 1 to 3,
# which cascades into recalculated expression/specificity statistics
# for all three data rows (2-4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3
        G = 1.030683333333333
        H = 3.09205
        K = 3
        M = 9.210619666666666
        N = 27.631859
        O = 0.133636377806767
        P = 0.133636377806767
        Q = 9.493232180105556
        R = 85.43908962095
        S = 0.133636377806767
        T = 0.133636377806767
    }
    3 = @{
        E = 3
        G = 1.030683333333333
        H = 3.09205
        K = 3
        M = 43.97212233333332
        N = 131.916367
        O = 0.6379891218794987
        P = 0.6379891218794989
        Q = 45.32133362026111
        R = 407.89200258235
        S = 0.6379891218794987
        T = 0.6379891218794989
    }
    4 = @{
        E = 3
        G = 1.030683333333333
        H = 3.09205
        K = 3
        M = 15.740255
        N = 47.220765
        O = 0.2283745003137342
        P = 0.2283745003137342
        Q = 16.22321849091667
        R = 146.00896641825
        S = 0.2283745003137342
        T = 0.2283745003137342
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
